$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Semana 16 de 2025 data refresh ---
# The refreshed extract has two additional event rows (310, 591), growing the
# table from 31 rows (A1:E31) to 33 rows (A1:E33). Insert blank rows at the right
# spots first (shifting everything below them down), then rewrite every data cell.

# New row for event "310" lands at (final) row 9, pushing the former row 9 ("330")
# and everything after it down by one.
$ws.Rows(9).Insert()

# New row for event "591" lands at (final) row 26 -- at this point in the sheet
# (after the first insert) that is still occupied by the former "620" row, so
# inserting here pushes it (and everything after) down by one more.
$ws.Rows(26).Insert()

# Make sure the two freshly inserted rows start out completely empty.
$ws.Rows(9).ClearContents()
$ws.Rows(26).ClearContents()

# Column A holds "evento" codes that look numeric (113, 115, ...); format it as
# text so they are stored the same way as in the source extract (t="inlineStr"/"s"),
# not auto-converted to numbers. (Row 1 is the "evento" text header, not a code,
# so it is left with its original General format.)
$ws.Range("A2:A33").NumberFormat = "@"

# Row 1
$ws.Range("A1").Value = "evento"
$ws.Range("B1").Value = "nom_eve"
$ws.Range("C1").Value = "Esperado"
$ws.Range("D1").Value = "Observado"
$ws.Range("E1").Value = "valor p"

# Row 2
$ws.Range("A2").Value = "113"
$ws.Range("B2").Value = "Desnutrici”n aguda en menores de 5 anos"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.37

# Row 3
$ws.Range("A3").Value = "115"
$ws.Range("B3").Value = "Cancer en menores de 18 anos"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1

# Row 4
$ws.Range("A4").Value = "155"
$ws.Range("B4").Value = "Cancer de la mama y cuello uterino"
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 0.15

# Row 5
$ws.Range("A5").Value = "210"
$ws.Range("B5").Value = "Dengue"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 13
$ws.Range("E5").Value = 0

# Row 6
$ws.Range("A6").Value = "215"
$ws.Range("B6").Value = "Defectos congenitos"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.37

# Row 7
$ws.Range("A7").Value = "220"
$ws.Range("B7").Value = "Dengue grave"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1

# Row 8
$ws.Range("A8").Value = "300"
$ws.Range("B8").Value = "Agresiones por animales potencialmente transmisores de rabia"
$ws.Range("C8").Value = 38
$ws.Range("D8").Value = 38
$ws.Range("E8").Value = 0.06

# Row 9
$ws.Range("A9").Value = "310"
$ws.Range("B9").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("D9").Value = 3
$ws.Range("E9").ClearContents()

# Row 10
$ws.Range("A10").Value = "330"
$ws.Range("B10").Value = "Hepatitis a"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.37

# Row 11
$ws.Range("A11").Value = "340"
$ws.Range("B11").Value = "Hepatitis b, c y coinfeccion hepatitis b y delta"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 0

# Row 12
$ws.Range("A12").Value = "342"
$ws.Range("B12").Value = "Enfermedades huerfanas - raras"
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 0.15

# Row 13
$ws.Range("A13").Value = "346"
$ws.Range("B13").Value = "Ira por virus nuevo"
$ws.Range("C13").Value = 13
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

# Row 14
$ws.Range("A14").Value = "348"
$ws.Range("B14").Value = "Infeccion respiratoria aguda grave irag inusitada"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 1

# Row 15
$ws.Range("A15").Value = "352"
$ws.Range("B15").Value = "Infecciones de sitio quirurgico asociadas a procedimiento medico quirurgico"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0.37

# Row 16
$ws.Range("A16").Value = "355"
$ws.Range("B16").Value = "Enfermedad transmitida por alimentos o agua (eta)"
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 1

# Row 17
$ws.Range("A17").Value = "356"
$ws.Range("B17").Value = "Intento de suicidio"
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = 0.07000000000000001

# Row 18
$ws.Range("A18").Value = "357"
$ws.Range("B18").Value = "Iad - infecciones asociadas a dispositivos - individual"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0.37

# Row 19
$ws.Range("A19").Value = "365"
$ws.Range("B19").Value = "Intoxicaciones"
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 0.05

# Row 20
$ws.Range("A20").Value = "455"
$ws.Range("B20").Value = "Leptospirosis"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0.37

# Row 21
$ws.Range("A21").Value = "465"
$ws.Range("B21").Value = "Malaria"
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 1

# Row 22
$ws.Range("A22").Value = "535"
$ws.Range("B22").Value = "Meningitis bacteriana y enfermedad meningoc”cica"
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 1

# Row 23
$ws.Range("A23").Value = "549"
$ws.Range("B23").Value = "Morbilidad materna extrema"
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 0.18

# Row 24
$ws.Range("A24").Value = "560"
$ws.Range("B24").Value = "Mortalidad perinatal y neonatal tardia"
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 0

# Row 25
$ws.Range("A25").Value = "580"
$ws.Range("B25").Value = "Mortalidad por dengue"
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 1

# Row 26
$ws.Range("A26").Value = "591"
$ws.Range("B26").Value = "Vigilancia integrada de muertes en menores de cinco anos por infeccion respiratoria aguda - enfermedad diarreica aguda y/o desnutricion"
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 1

# Row 27
$ws.Range("A27").Value = "620"
$ws.Range("B27").Value = "Parotiditis"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0.37

# Row 28
$ws.Range("A28").Value = "740"
$ws.Range("B28").Value = "Sifilis congenita"
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 1

# Row 29
$ws.Range("A29").Value = "750"
$ws.Range("B29").Value = "Sifilis gestacional"
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0.37

# Row 30
$ws.Range("A30").Value = "800"
$ws.Range("B30").Value = "Tos ferina"
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 1

# Row 31
$ws.Range("A31").Value = "813"
$ws.Range("B31").Value = "Tuberculosis"
$ws.Range("C31").Value = 8
$ws.Range("D31").Value = 8
$ws.Range("E31").Value = 0.14

# Row 32
$ws.Range("A32").Value = "831"
$ws.Range("B32").Value = "Varicela individual"
$ws.Range("C32").Value = 11
$ws.Range("D32").Value = 4
$ws.Range("E32").Value = 0.01

# Row 33
$ws.Range("A33").Value = "850"
$ws.Range("B33").Value = "Vih/sida/mortalidad por sida"
$ws.Range("C33").Value = 7
$ws.Range("D33").Value = 6
$ws.Range("E33").Value = 0.15
